# Reliability Measurements workbook update:
#  - rename the existing "Sheet1" to "9-2-14"
#  - duplicate it (preserving layout/styles/formulas) as "9-3-14"
#  - trim the new sheet down to a single (1mL Syringe) trial block
#  - replace that block's data with the "9-3-14" trial results (5 trials)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the second day's sheet as a copy of the first ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "9-2-14"
$ws2.Name = "9-3-14"

# Row-height touch-up on 9-2-14 (auto height settles to single-line values)
$ws1.Rows.Item(1).RowHeight = 15
$ws1.Rows.Item(14).RowHeight = 15.75
$ws1.Rows.Item(29).RowHeight = 15.75

# --- 9-3-14: keep only the first ("Test Pipette (1mL Syringe)") block ---
$ws2.Rows("9:38").Delete()
# drop the 6th trial row - the new data set only has 5 trials
$ws2.Rows("8:8").Delete()

$ws2.Rows.Item(2).RowHeight = 75.75

# New trial measurements recorded on 9-3-14
$trialNum = @(1, 2, 3, 4, 5)
$dispensed = @(0.2028, 0.2044, 0.2014, 0.2004, 0.201)
$syringe = @(0.1994, 0.2001, 0.1997, 0.1995, 0.1962)

for ($i = 0; $i -lt 5; $i++) {
    $r = 3 + $i
    $ws2.Cells.Item($r, 1).Value = $trialNum[$i]
    $ws2.Cells.Item($r, 2).Value = $dispensed[$i]
    $ws2.Cells.Item($r, 3).Value = $syringe[$i]
}

# Recreate the "Average" summary row (formatting copied down from the last trial row)
$ws2.Range("A7:C7").Copy()
$ws2.Range("A8:C8").PasteSpecial(-4122)
$ws2.Rows.Item(8).RowHeight = 15.75

$ws2.Cells.Item(8, 1).Value = "Average "
$ws2.Range("B8").Formula = "=(B3+B4+B5+B6+B7)/5"
$ws2.Range("C8").Formula = "=(C3+C5+C4+C6+C7)/5"

# column widths trimmed slightly on the new sheet
$ws2.Columns.Item(2).ColumnWidth = 43.57
$ws2.Columns.Item(3).ColumnWidth = 58.43

$ws2.Range("B11").Select()

# leave the first day's sheet as the active tab
$ws1.Activate()
$ws1.Range("G12").Select()

Write-Output "done"
